# Applies the "output generated at 456a3b4" update to 上海-漫展信息.xlsx
# - Sheet 展览 (exhibitions): refresh "want-to-go" counts (col F)
# - Sheet 演出 (performances): refresh "want-to-go" counts (col F)
# - Sheet 本地生活 (local life): refresh "want-to-go" counts (col F) and
#   insert a brand-new row for the 银魂快闪店 pop-up (pushing the 柯南展
#   row from 13 -> 14, with its own F value bumped too)
# - Sheet 全部类型 (all types, a merged/sorted view of the above): refresh
#   the corresponding "want-to-go" counts (col F)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 241
$ws.Range("F4").Value = 503
$ws.Range("F5").Value = 2181
$ws.Range("F7").Value = 8126
$ws.Range("F8").Value = 87
$ws.Range("F10").Value = 260
$ws.Range("F11").Value = 1787
$ws.Range("F12").Value = 1574
$ws.Range("F13").Value = 1321
$ws.Range("F14").Value = 198
$ws.Range("F15").Value = 4254
$ws.Range("F16").Value = 6145
$ws.Range("F17").Value = 733
$ws.Range("F18").Value = 50
$ws.Range("F19").Value = 1160
$ws.Range("F20").Value = 1246
$ws.Range("F21").Value = 452
$ws.Range("F22").Value = 6338
$ws.Range("F25").Value = 4293
$ws.Range("F26").Value = 278
$ws.Range("F27").Value = 710
$ws.Range("F28").Value = 1975
$ws.Range("F29").Value = 1178
$ws.Range("F30").Value = 323
$ws.Range("F31").Value = 1050
$ws.Range("F32").Value = 42
$ws.Range("F33").Value = 37
$ws.Range("F34").Value = 64
$ws.Range("F36").Value = 1173
$ws.Range("F38").Value = 1891
$ws.Range("F40").Value = 424
$ws.Range("F41").Value = 164
$ws.Range("F42").Value = 1184
$ws.Range("F44").Value = 66
$ws.Range("F45").Value = 1092
$ws.Range("F47").Value = 80
$ws.Range("F48").Value = 182

# ---------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 81
$ws.Range("F4").Value = 995
$ws.Range("F9").Value = 3
$ws.Range("F13").Value = 688
$ws.Range("F14").Value = 386
$ws.Range("F23").Value = 151
$ws.Range("F24").Value = 100
$ws.Range("F32").Value = 275
$ws.Range("F33").Value = 46

# ---------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 461
$ws.Range("F6").Value = 1582
$ws.Range("F7").Value = 482
$ws.Range("F8").Value = 3124
$ws.Range("F9").Value = 1003
$ws.Range("F10").Value = 1125
$ws.Range("F11").Value = 1352
$ws.Range("F12").Value = 1633

# Insert a fresh row 13 (everything currently on row 13, the Conan
# exhibit, slides down to row 14).
$ws.Rows("13:13").Insert()

# Give the new row's cells the same look as the surrounding data rows
# (copy the format down from row 12) before filling in values, so the
# new cells don't pick up an odd inherited style.
$ws.Range("A12:I12").Copy()
$ws.Range("A13:I13").PasteSpecial(-4122)

$ws.Range("A13").Value = 12

# Columns that look like dates (B, E) get auto-converted to real Excel
# dates on plain assignment; force them to text first so they land as
# literal strings like the rest of the sheet, then strip the leftover
# "@" number format (copied again from a plain neighbour) without
# disturbing the text that's already been committed to the cell.
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2024-08-03"
$ws.Range("C13").Value = "上海· 银魂快闪店in大悦城 游艺体验套装"
$ws.Range("D13").Value = "西藏北路166号（地铁8号线曲阜路下） 静安大悦城"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2024.08.03 00:00-09.01 23:59"
$ws.Range("F13").Value = 40
$ws.Range("G13").Value = 68
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=89345"
$ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202407/VvQSpDGa1721117043746.png"

$ws.Range("B12:E12").Copy()
$ws.Range("B13:E13").PasteSpecial(-4122)

# Old row 13 (now row 14): bump its index number to account for the
# newly-inserted row above it, and refresh its want-to-go count.
$ws.Range("A14").Value = 13
$ws.Range("F14").Value = 52

# ---------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 461
$ws.Range("F3").Value = 241
$ws.Range("F5").Value = 503
$ws.Range("F6").Value = 482
$ws.Range("F7").Value = 3124
$ws.Range("F8").Value = 2181
$ws.Range("F9").Value = 87
$ws.Range("F10").Value = 1003
$ws.Range("F12").Value = 3
$ws.Range("F14").Value = 260
$ws.Range("F15").Value = 1787
$ws.Range("F16").Value = 1574
$ws.Range("F17").Value = 1352
$ws.Range("F18").Value = 688
$ws.Range("F19").Value = 198
$ws.Range("F20").Value = 1633
$ws.Range("F21").Value = 4254
$ws.Range("F22").Value = 386
$ws.Range("F24").Value = 733
$ws.Range("F25").Value = 50
$ws.Range("F26").Value = 1160
$ws.Range("F27").Value = 1246
$ws.Range("F28").Value = 452
$ws.Range("F29").Value = 6339
$ws.Range("F31").Value = 710
$ws.Range("F32").Value = 1975
$ws.Range("F33").Value = 1178
$ws.Range("F34").Value = 323
$ws.Range("F35").Value = 37
$ws.Range("F36").Value = 64
$ws.Range("F37").Value = 100
$ws.Range("F40").Value = 1891
$ws.Range("F42").Value = 424
$ws.Range("F43").Value = 1184
$ws.Range("F46").Value = 275
$ws.Range("F47").Value = 1092
$ws.Range("F48").Value = 182

Write-Host "edit applied"
